$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A70").Value = 41384
$ws.Range("B70").Value = "Finished execution model chapter, wrote memory chapter, added sample file, FINISHED OPENCL chapter"

$ws.Range("B71").Select()
